# Commit: "added 4wk low sales check"
# Re-ran the forecast after adding a 4-week low-sales check, which shifted the
# MyForecast / Inventory Coverage / Seasonality Index numbers (and a couple of
# derived risk labels) on the "Forecast Comparison" sheet, and the roll-up
# totals/extremes on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet ---

# Row 2 (W10)
$wsForecast.Range("D2").Value = 162
$wsForecast.Range("H2").Value = 9.93
$wsForecast.Range("L2").Value = 1.18

# Row 3 (W11)
$wsForecast.Range("D3").Value = 173
$wsForecast.Range("H3").Value = 8.359999999999999
$wsForecast.Range("L3").Value = 1.04

# Row 4 (W12)
$wsForecast.Range("D4").Value = 172
$wsForecast.Range("H4").Value = 7.41
$wsForecast.Range("L4").Value = 1.16

# Row 5 (W13)
$wsForecast.Range("D5").Value = 155
$wsForecast.Range("H5").Value = 7.11
$wsForecast.Range("L5").Value = 1.1

# Row 6 (W14)
$wsForecast.Range("D6").Value = 136
$wsForecast.Range("H6").Value = 6.96
$wsForecast.Range("L6").Value = 1

# Row 7 (W15)
$wsForecast.Range("D7").Value = 126
$wsForecast.Range("H7").Value = 6.44
$wsForecast.Range("L7").Value = 1.18

# Row 8 (W16)
$wsForecast.Range("D8").Value = 130
$wsForecast.Range("H8").Value = 5.27
$wsForecast.Range("L8").Value = 1.15

# Row 9 (W17)
$wsForecast.Range("D9").Value = 137
$wsForecast.Range("H9").Value = 4.05
$wsForecast.Range("L9").Value = 0.91

# Row 10 (W18)
$wsForecast.Range("D10").Value = 135
$wsForecast.Range("H10").Value = 3.1
$wsForecast.Range("L10").Value = 1

# Row 11 (W19)
$wsForecast.Range("D11").Value = 119
$wsForecast.Range("H11").Value = 2.38
$wsForecast.Range("L11").Value = 1.19

# Row 12 (W20)
$wsForecast.Range("D12").Value = 102
$wsForecast.Range("H12").Value = 1.61
$wsForecast.Range("J12").Value = "Normal"
$wsForecast.Range("L12").Value = 1.12

# Row 13 (W21)
$wsForecast.Range("D13").Value = 98
$wsForecast.Range("H13").Value = 0.63
$wsForecast.Range("I13").Value = "Low"
$wsForecast.Range("L13").Value = 0.85

# Row 14 (W22)
$wsForecast.Range("D14").Value = 107
$wsForecast.Range("L14").Value = 0.86

# Row 15 (W23)
$wsForecast.Range("D15").Value = 114
$wsForecast.Range("L15").Value = 1.08

# Row 16 (W24)
$wsForecast.Range("D16").Value = 108
$wsForecast.Range("L16").Value = 1.18

# Row 17 (W25)
$wsForecast.Range("D17").Value = 90
$wsForecast.Range("L17").Value = 0.88

# --- "Summary" sheet ---
# These values are stored as text (not numbers) in the sheet, so a leading
# apostrophe is used to force text entry, same as typing them in Excel.
$wsSummary.Range("B9").Value  = "'2064"   # Total Forecast (16 Weeks)
$wsSummary.Range("B10").Value = "'1191"   # Total Forecast (8 Weeks)
$wsSummary.Range("B11").Value = "'662"    # Total Forecast (4 Weeks)
$wsSummary.Range("B12").Value = "'173"    # Max Forecast
$wsSummary.Range("B14").Value = "'90"     # Min Forecast
